$d = $word.ActiveDocument

# --- Change 1: Title paragraph - merge "TRACKTOR v0.1 " run and "user manual" run ---
# The title paragraph originally consists of two runs:
#   Run 1: "TRACKTOR v0.1 " (with a trailing space, xml:space="preserve")
#   Run 2: "user manual"
# The target document merges them into a single run containing the full text.
$prefix = "TRACKTOR v0.1 " + [char]0x2013 + " "
$prefixLen = $prefix.Length
$fullTitle = "TRACKTOR v0.1 " + [char]0x2013 + " user manual"
$fullLen = $fullTitle.Length

$titleRun2 = $d.Range($prefixLen, $fullLen)
if ($titleRun2.Text -eq "user manual") {
    $titleRun2.Delete()
    $titleRun1 = $d.Range(0, $prefixLen)
    $titleRun1.InsertAfter("user manual")
}

# --- Change 2: pip3 -> pip (install instructions) ---
$d.Content.Find.Execute("pip3 install git+https://github.com/vivekhsridhar/tracktor.git",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pip install git+https://github.com/vivekhsridhar/tracktor.git", 2) | Out-Null
